$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# The sheet contains an account-entries import sample (rows 2-18). A new
# entry (duplicate of the former last row - Investor 4 / Folio 10 /
# Investable Capital Percentage / 0.06) was inserted at row 2, shifting all
# the other rows down by one; the table still ends at row 18 because the
# previously-last row's data now simply lives one row earlier than before.
# We therefore rewrite the value of every data cell (rows 2-18) to its new
# "after" content, and fix up the two rows whose cell styling changed
# (row 2 lost the A/B/C styling that a normal data row has, and row 18
# gained it).
# --------------------------------------------------------------------------

$ws.Cells.Item(2,1).Value = "Investor 4"
$ws.Cells.Item(2,2).Value = "SAAS Fund"
$ws.Cells.Item(2,3).Value = 10
$ws.Cells.Item(2,4).Value = 44961
$ws.Cells.Item(2,5).Value = "Fees"
$ws.Cells.Item(2,6).Value = "Investable Capital Percentage"
$ws.Cells.Item(2,7).Value = 0.06
$ws.Cells.Item(2,8).ClearContents()
$ws.Cells.Item(2,9).Value = "Pool"

$ws.Cells.Item(3,1).Value = "Investor 1"
$ws.Cells.Item(3,2).Value = "SAAS Fund"
$ws.Cells.Item(3,3).Value = 6
$ws.Cells.Item(3,4).Value = 44960
$ws.Cells.Item(3,5).Value = "Expense"
$ws.Cells.Item(3,6).Value = "Management Fees"
$ws.Cells.Item(3,7).Value = 1000
$ws.Cells.Item(3,8).ClearContents()
$ws.Cells.Item(3,9).Value = "Pool"

$ws.Cells.Item(4,1).Value = "Investor 2"
$ws.Cells.Item(4,2).Value = "SAAS Fund"
$ws.Cells.Item(4,3).Value = 7
$ws.Cells.Item(4,4).Value = 44960
$ws.Cells.Item(4,5).Value = "Interest"
$ws.Cells.Item(4,6).Value = "Accrued Interest"
$ws.Cells.Item(4,7).Value = 2000
$ws.Cells.Item(4,8).Value = "From FDs"
$ws.Cells.Item(4,9).Value = "Pool"

$ws.Cells.Item(5,1).Value = "Investor 2"
$ws.Cells.Item(5,2).Value = "SAAS Fund"
$ws.Cells.Item(5,3).Value = 7
$ws.Cells.Item(5,4).Value = 44960
$ws.Cells.Item(5,5).Value = "Expense"
$ws.Cells.Item(5,6).Value = "Management Fees"
$ws.Cells.Item(5,7).Value = 1000
$ws.Cells.Item(5,8).ClearContents()
$ws.Cells.Item(5,9).Value = "Pool"

$ws.Cells.Item(6,1).Value = "Investor 3"
$ws.Cells.Item(6,2).Value = "SAAS Fund"
$ws.Cells.Item(6,3).Value = 8
$ws.Cells.Item(6,4).Value = 44961
$ws.Cells.Item(6,5).Value = "Expense"
$ws.Cells.Item(6,6).Value = "Management Fees"
$ws.Cells.Item(6,7).Value = 1000
$ws.Cells.Item(6,8).ClearContents()
$ws.Cells.Item(6,9).Value = "Pool"

$ws.Cells.Item(7,1).Value = "Investor 4"
$ws.Cells.Item(7,2).Value = "SAAS Fund"
$ws.Cells.Item(7,3).Value = 9
$ws.Cells.Item(7,4).Value = 44961
$ws.Cells.Item(7,5).Value = "Interest"
$ws.Cells.Item(7,6).Value = "Accrued Interest"
$ws.Cells.Item(7,7).Value = 2000
$ws.Cells.Item(7,8).Value = "From FDs"
$ws.Cells.Item(7,9).Value = "Pool"

$ws.Cells.Item(8,1).Value = "Investor 4"
$ws.Cells.Item(8,2).Value = "SAAS Fund"
$ws.Cells.Item(8,3).Value = 9
$ws.Cells.Item(8,4).Value = 44961
$ws.Cells.Item(8,5).Value = "Expense"
$ws.Cells.Item(8,6).Value = "Management Fees"
$ws.Cells.Item(8,7).Value = 1000
$ws.Cells.Item(8,8).ClearContents()
$ws.Cells.Item(8,9).Value = "Pool"

$ws.Cells.Item(9,1).Value = "Investor 4"
$ws.Cells.Item(9,2).Value = "SAAS Fund"
$ws.Cells.Item(9,3).Value = 10
$ws.Cells.Item(9,4).Value = 44961
$ws.Cells.Item(9,5).Value = "Expense"
$ws.Cells.Item(9,6).Value = "Management Fees"
$ws.Cells.Item(9,7).Value = 1000
$ws.Cells.Item(9,8).ClearContents()
$ws.Cells.Item(9,9).Value = "Pool"

$ws.Cells.Item(10,1).Value = "Investor 1"
$ws.Cells.Item(10,2).Value = "SAAS Fund"
$ws.Cells.Item(10,3).Value = 6
$ws.Cells.Item(10,4).Value = 44960
$ws.Cells.Item(10,5).Value = "Fees"
$ws.Cells.Item(10,6).Value = "Setup Fees"
$ws.Cells.Item(10,7).Value = 500
$ws.Cells.Item(10,8).ClearContents()
$ws.Cells.Item(10,9).Value = "Pool"

$ws.Cells.Item(11,1).Value = "Investor 2"
$ws.Cells.Item(11,2).Value = "SAAS Fund"
$ws.Cells.Item(11,3).Value = 7
$ws.Cells.Item(11,4).Value = 44960
$ws.Cells.Item(11,5).Value = "Fees"
$ws.Cells.Item(11,6).Value = "Setup Fees"
$ws.Cells.Item(11,7).Value = 500
$ws.Cells.Item(11,8).ClearContents()
$ws.Cells.Item(11,9).Value = "Pool"

$ws.Cells.Item(12,1).Value = "Investor 3"
$ws.Cells.Item(12,2).Value = "SAAS Fund"
$ws.Cells.Item(12,3).Value = 8
$ws.Cells.Item(12,4).Value = 44961
$ws.Cells.Item(12,5).Value = "Fees"
$ws.Cells.Item(12,6).Value = "Setup Fees"
$ws.Cells.Item(12,7).Value = 500
$ws.Cells.Item(12,8).ClearContents()
$ws.Cells.Item(12,9).Value = "Pool"

$ws.Cells.Item(13,1).Value = "Investor 4"
$ws.Cells.Item(13,2).Value = "SAAS Fund"
$ws.Cells.Item(13,3).Value = 9
$ws.Cells.Item(13,4).Value = 44961
$ws.Cells.Item(13,5).Value = "Fees"
$ws.Cells.Item(13,6).Value = "Setup Fees"
$ws.Cells.Item(13,7).Value = 500
$ws.Cells.Item(13,8).ClearContents()
$ws.Cells.Item(13,9).Value = "Pool"

$ws.Cells.Item(14,1).Value = "Investor 4"
$ws.Cells.Item(14,2).Value = "SAAS Fund"
$ws.Cells.Item(14,3).Value = 10
$ws.Cells.Item(14,4).Value = 44961
$ws.Cells.Item(14,5).Value = "Fees"
$ws.Cells.Item(14,6).Value = "Setup Fees"
$ws.Cells.Item(14,7).Value = 500
$ws.Cells.Item(14,8).ClearContents()
$ws.Cells.Item(14,9).Value = "Pool"

$ws.Cells.Item(15,1).Value = "Investor 1"
$ws.Cells.Item(15,2).Value = "SAAS Fund"
$ws.Cells.Item(15,3).Value = 6
$ws.Cells.Item(15,4).Value = 44960
$ws.Cells.Item(15,5).Value = "Fees"
$ws.Cells.Item(15,6).Value = "Investable Capital Percentage"
$ws.Cells.Item(15,7).Value = 0.05
$ws.Cells.Item(15,8).ClearContents()
$ws.Cells.Item(15,9).Value = "Pool"

$ws.Cells.Item(16,1).Value = "Investor 2"
$ws.Cells.Item(16,2).Value = "SAAS Fund"
$ws.Cells.Item(16,3).Value = 7
$ws.Cells.Item(16,4).Value = 44960
$ws.Cells.Item(16,5).Value = "Fees"
$ws.Cells.Item(16,6).Value = "Investable Capital Percentage"
$ws.Cells.Item(16,7).Value = 0.06
$ws.Cells.Item(16,8).ClearContents()
$ws.Cells.Item(16,9).Value = "Pool"

$ws.Cells.Item(17,1).Value = "Investor 3"
$ws.Cells.Item(17,2).Value = "SAAS Fund"
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 44961
$ws.Cells.Item(17,5).Value = "Fees"
$ws.Cells.Item(17,6).Value = "Investable Capital Percentage"
$ws.Cells.Item(17,7).Value = 0.05
$ws.Cells.Item(17,8).ClearContents()
$ws.Cells.Item(17,9).Value = "Pool"

$ws.Cells.Item(18,1).Value = "Investor 4"
$ws.Cells.Item(18,2).Value = "SAAS Fund"
$ws.Cells.Item(18,3).Value = 9
$ws.Cells.Item(18,4).Value = 44961
$ws.Cells.Item(18,5).Value = "Fees"
$ws.Cells.Item(18,6).Value = "Investable Capital Percentage"
$ws.Cells.Item(18,7).Value = 0.06
$ws.Cells.Item(18,8).ClearContents()
$ws.Cells.Item(18,9).Value = "Pool"


# Row 2 is the newly inserted row: it carries no explicit per-cell style on
# columns A/B/C (same as the old trailing "scratch" row used to have).
$ws.Cells.Item(2,1).Style = "Normal"
$ws.Cells.Item(2,2).Style = "Normal"
$ws.Cells.Item(2,3).Style = "Normal"

# Row 18 now holds a regular data row, so it picks up the normal styling
# used by every other data row in the table.
$ws.Cells.Item(18,1).Style = "Normal 3"
$ws.Cells.Item(18,2).Style = "Normal 5"
$ws.Cells.Item(18,3).Style = "Normal 3"


# Reflect the new selection left behind in the sheet view (entire row 2 was
# selected when the row was inserted).
$ws.Range("A2:XFD2").Select()
